$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "bigger_dem_raster"
$ws.Range("B10").Value = "data/dtm_big_area_depth_padded.tif"
$ws.Range("C10").Value = "Necessary when the mesh is larger than the DEM and has some cell centers outside of it. If not provided, such mesh would get NaN values!"

$ws.Range("C10").Select()
